$p = $ppt.ActivePresentation
try {
  $newDesign = $p.Designs.Add()
  Write-Output ("Added, Count=" + $p.Designs.Count)
} catch {
  Write-Output ("ERR: " + $_.Exception.Message)
}
